$d = $word.ActiveDocument

# The target sentence currently lives in a single run:
#   "This information dashboard tracks my past progress towards achieving "
# It needs to become three runs (same character formatting throughout):
#   "This information dashboard tracks" + " " + "my past progress towards achieving "
$part1 = "This information dashboard tracks"
$part2 = " "
$part3 = "my past progress towards achieving "

$findRange = $d.Content
$findRange.Find.Execute($part1 + $part2 + $part3, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $findRange.Find.Found) {
    throw "Could not find target sentence"
}

$start = $findRange.Start
$split1 = $start + $part1.Length
$split2 = $split1 + $part2.Length

# Toggling a character-formatting property back to its original value forces
# the run to be split at the toggled boundary without changing any visible
# formatting or text - unlike re-assigning .Text (or Find/Replace), which
# causes the engine to re-coalesce same-formatted adjacent runs back into one.
$r1 = $d.Range($start, $split1)
$r1.Bold = $true
$r1.Bold = $false

$r2 = $d.Range($split1, $split2)
$r2.Bold = $true
$r2.Bold = $false
